$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-31 07:55:01"
$wsZhCn.Range("E3").Value = "2016-03-31 07:55:01"
$wsZhCn.Range("H2").Value = "2016-03-31 07:56:05"
$wsZhCn.Range("H3").Value = "2016-03-31 07:56:05"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-31 07:55:13"
$wsDeDe.Range("E3").Value = "2016-03-31 07:55:13"
$wsDeDe.Range("H2").Value = "2016-03-31 07:56:22"
$wsDeDe.Range("H3").Value = "2016-03-31 07:56:22"
